# Update the "mapa_interactivo" workbook.
#
# The record Caso 6002 (LA PLATA AV. 832) is reassigned from provider
# "PEBCOM" to provider "NEW". Concretely:
#   1. On the "General" master sheet, its "Proveedor Asignado" cell (F324)
#      changes from PEBCOM to NEW.
#   2. On the "PEBCOM" sheet, the row holding this record (row 73) is
#      removed, and every following row shifts up by one (dimension
#      shrinks from P92 to P91).
#   3. On the "NEW" sheet, a new row is inserted at the top of the data
#      block (row 55) holding this same record (now flagged Proveedor =
#      NEW), and every following row shifts down by one (dimension grows
#      from P76 to P77).

$wb = $excel.ActiveWorkbook

# --- 1. Update the General sheet ---------------------------------------
$general = $wb.Worksheets.Item("General")
$general.Range("F324").Value = "NEW"

# --- 2 & 3. Move the row from PEBCOM to the top of the NEW block -------
$pebcom = $wb.Worksheets.Item("PEBCOM")
$new = $wb.Worksheets.Item("NEW")

# Insert a blank row at the destination, then copy the source row's
# cells into it (preserves text/number typing and avoids Excel
# reinterpreting date-like text such as "8/12/2025" as a date serial).
$new.Rows.Item(55).Insert()
$pebcom.Rows.Item(73).Copy()
$new.Rows.Item(55).PasteSpecial(-4104)
$new.Range("F55").Value = "NEW"

# Now remove the original row from PEBCOM.
$pebcom.Rows.Item(73).Delete()
